# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 12:04"

# Row 4 (Estados Unidos) updated counts
$ws.Range("B4").Value = 1322164
$ws.Range("C4").Value = 379
$ws.Range("E4").Value = 1019799

# Row 37 updated counts
$ws.Range("E37").Value = 7462
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 926

# Row 55 updated counts
$ws.Range("B55").Value = 5880
$ws.Range("C55").Value = 142
$ws.Range("E55").Value = 1615
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 265

# Rows 212/213: "Butan" and "Islas Virgenes Britanicas" swap places
# (country names swap, along with their Recuperados (D) / Muertes (H) values)
$ws.Range("A212").Value = "Islas Virgenes Britanicas"
$ws.Range("D212").Value = 4
$ws.Range("H212").Value = 1

$ws.Range("A213").Value = "Butan"
$ws.Range("D213").Value = 5
$ws.Range("H213").Value = 0
